$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2 currently holds the numeric value 56122.
# Convert it in place to a text string "56122" (same displayed text, but stored as a string)
# without leaving a lingering custom style on the cell.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "56122"
$ws.Range("A2").ClearFormats()

# Add the new row: A3 = 79086 (numeric)
$ws.Range("A3").Value = 79086
